# CDS test cases for Experimental Strategy and File Type filters
# Rework the "startup" sheet from Cases/Samples/Files tabs (filtered on
# study.experimental_strategy_and_data_subtypes) to Participants/Samples/Files
# tabs filtered on file.experimental_strategy_and_data_subtypes, each query
# now joining through (samp)<--(f:file) and capped with LIMIT 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- StatQuery (shared by column C on every row) ------------------------
$statQuery = @"
MATCH (s:study)<--(p:participant)
OPTIONAL MATCH (p)<--(samp:sample)
MATCH (samp)<--(f:file)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH DISTINCT samp,diag,s,p,f
WHERE f.experimental_strategy_and_data_subtypes in ["Methylation Array"]
RETURN
    count(distinct s) AS Studies,
    count(distinct p) AS Participants,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Files``
"@

# ---- Participants tab query (row 2) --------------------------------------
$participantsQuery = @"
MATCH (s:study)<--(p:participant)
OPTIONAL MATCH (p)<--(samp:sample)
MATCH (samp)<--(f:file)
WHERE f.experimental_strategy_and_data_subtypes in ["Methylation Array"]
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN   
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
 coalesce(p.gender,'') as ``Gender``,
 coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER By p.participant_id LIMIT 100
"@

# ---- Samples tab query (row 3) -------------------------------------------
$samplesQuery = @"
MATCH (s:study)<--(p:participant)<--(samp:sample)
MATCH (samp)<--(f:file)
WHERE f.experimental_strategy_and_data_subtypes in ["Methylation Array"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as ``Sample ID``,
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
 coalesce(samp.sample_tumor_status,'') as ``Tumor``,
coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER By samp.sample_id LIMIT 100
"@

# ---- Files tab query (row 4) ----------------------------------------------
$filesQuery = @"
MATCH (s:study)<--(p:participant)<--(samp:sample)
MATCH (samp)<--(f:file)
WHERE f.experimental_strategy_and_data_subtypes in ["Methylation Array"]
WITH f,p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER By f.file_name LIMIT 100
"@

$neo4jDataFile = "TC02_CDS_Filter_ExprStrtgies-MethylationArray_Neo4jData.xlsx"
$webDataFile = "TC02_CDS_Filter_ExprStrtgies-MethylationArray_WebData.xlsx"

# Row 2: was CasesTab -> now ParticipantsTab
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jDataFile
$ws.Range("E2").Value = $webDataFile

# Row 3: SamplesTab (unchanged name, new query)
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jDataFile
$ws.Range("E3").Value = $webDataFile

# Row 4: FilesTab (unchanged name, new query)
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jDataFile
$ws.Range("E4").Value = $webDataFile

# ---- Row heights: all three data rows now render at the same height -----
$ws.Rows(2).RowHeight = 186
$ws.Rows(3).RowHeight = 186
$ws.Rows(4).RowHeight = 186

# ---- Column widths: nudge to match the refreshed autofit layout ---------
$ws.Columns(1).ColumnWidth = 11.4167
$ws.Columns(2).ColumnWidth = 74.7556
$ws.Columns(3).ColumnWidth = 74.7556
$ws.Columns(4).ColumnWidth = 69.4167
$ws.Columns(5).ColumnWidth = 62.5944

# ---- Selection moves from B2 to B5 ---------------------------------------
$ws.Range("B5").Select()
